$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Insert a new row at row 81 (shifts the existing rows 81-83 down to 82-84)
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row with the new sale record
$ws.Cells.Item(81, 1).Value2 = "HEMANTH HOME APPLIANCES"
$ws.Cells.Item(81, 2).Value2 = "kengeri satellite town"
$ws.Cells.Item(81, 3).Value2 = 145
$ws.Cells.Item(81, 4).Value2 = 45337
$ws.Cells.Item(81, 5).Value2 = 5000
$ws.Cells.Item(81, 6).Value2 = 0
$ws.Cells.Item(81, 7).Value2 = 0
$ws.Cells.Item(81, 8).Value2 = 0

# Match the bordered look of the rest of the table for the new row
$ws.Range("A81:H81").Borders.LineStyle = 1

# The Sales sheet's AutoFilter range now covers one additional row
$wb.Names.Item('Sales!_FilterDatabase').RefersTo = '=Sales!$A$1:$I$84'

# Leave the selection where the author left it after the edit
$ws.Range("B86").Select()
